$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4.946000000000002
$ws.Range("C4").Value = -14.2702
$ws.Range("D4").Value = -7.892199999999998

$ws.Range("C5").Value = -14.78610000000002

$ws.Range("B7").Value = 6.180499999999999

$ws.Range("C8").Value = -11.88649999999999

$ws.Range("D9").Value = -6.737999999999997

$ws.Range("B16").Value = 9.513500000000004
$ws.Range("C16").Value = -12.0769

$ws.Range("D18").Value = -8.457599999999996
